# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps on the per-locale report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 05:57:58"
$wsZhCn.Range("E3").Value = "2016-03-22 05:57:58"
$wsZhCn.Range("H2").Value = "2016-03-22 05:58:40"
$wsZhCn.Range("H3").Value = "2016-03-22 05:58:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 05:58:06"
$wsDeDe.Range("E3").Value = "2016-03-22 05:58:06"
$wsDeDe.Range("H2").Value = "2016-03-22 05:58:53"
$wsDeDe.Range("H3").Value = "2016-03-22 05:58:53"
